$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "template"
